$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = '24.918.53'
$ws.Range("E2").Value = '  +1.95%  '

$ws.Range("D3").Value = '1.708.35'
$ws.Range("E3").Value = '  +1.74%  '

$r = $ws.Range("D4")
$r.NumberFormat = "@"
$r.Value = '1.004'
$r.Style = "Normal"
$ws.Range("E4").Value = '  -0.13%  '

$r = $ws.Range("D5")
$r.NumberFormat = "@"
$r.Value = '313.43'
$r.Style = "Normal"
$ws.Range("E5").Value = '  +2.17%  '

$r = $ws.Range("D6")
$r.NumberFormat = "@"
$r.Value = '0.9987'
$r.Style = "Normal"
$ws.Range("E6").Value = '  -0.18%  '

$r = $ws.Range("D7")
$r.NumberFormat = "@"
$r.Value = '0.3746'
$r.Style = "Normal"
$ws.Range("E7").Value = '  +1.11%  '

$r = $ws.Range("D8")
$r.NumberFormat = "@"
$r.Value = '49.46'
$r.Style = "Normal"
$ws.Range("E8").Value = '  +3.51%  '

$r = $ws.Range("D9")
$r.NumberFormat = "@"
$r.Value = '0.3444'
$r.Style = "Normal"
$ws.Range("E9").Value = '  +0.51%  '

$r = $ws.Range("D10")
$r.NumberFormat = "@"
$r.Value = '1.224'
$r.Style = "Normal"
$ws.Range("E10").Value = '  +4.94%  '

$r = $ws.Range("D11")
$r.NumberFormat = "@"
$r.Value = '0.07551'
$r.Style = "Normal"
$ws.Range("E11").Value = '  +4.25%  '

$r = $ws.Range("D12")
$r.NumberFormat = "@"
$r.Value = '1.000'
$r.Style = "Normal"
$ws.Range("E12").Value = '  -0.53%  '

$r = $ws.Range("D13")
$r.NumberFormat = "@"
$r.Value = '21.25'
$r.Style = "Normal"
$ws.Range("E13").Value = '  +5.42%  '

$r = $ws.Range("D14")
$r.NumberFormat = "@"
$r.Value = '6.316'
$r.Style = "Normal"
$ws.Range("E14").Value = '  +3.72%  '

$r = $ws.Range("D15")
$r.NumberFormat = "@"
$r.Value = '7.077'
$r.Style = "Normal"
$ws.Range("E15").Value = '  +5.22%  '

$ws.Range("D16").Value = '1.707.95'
$ws.Range("E16").Value = '  +1.65%  '

$r = $ws.Range("D17")
$r.NumberFormat = "@"
$r.Value = '0.00001132'
$r.Style = "Normal"
$ws.Range("E17").Value = '  +2.56%  '

$r = $ws.Range("D18")
$r.NumberFormat = "@"
$r.Value = '0.06726'
$r.Style = "Normal"
$ws.Range("E18").Value = '  +0.71%  '

$ws.Range("E19").Value = '  -0.38%  '

$r = $ws.Range("D20")
$r.NumberFormat = "@"
$r.Value = '84.19'
$r.Style = "Normal"
$ws.Range("E20").Value = '  +3.95%  '

$r = $ws.Range("D21")
$r.NumberFormat = "@"
$r.Value = '17.31'
$r.Style = "Normal"
$ws.Range("E21").Value = '  +5.58%  '

$r = $ws.Range("D22")
$r.NumberFormat = "@"
$r.Value = '6.397'
$r.Style = "Normal"
$ws.Range("E22").Value = '  +4.92%  '

$r = $ws.Range("D23")
$r.NumberFormat = "@"
$r.Value = '13.08'
$r.Style = "Normal"
$ws.Range("E23").Value = '  +7.92%  '

$ws.Range("D24").Value = '24.907.90'
$ws.Range("E24").Value = '  +2.24%  '

$r = $ws.Range("D25")
$r.NumberFormat = "@"
$r.Value = '2.452'
$r.Style = "Normal"
$ws.Range("E25").Value = '  -0.64%  '

$r = $ws.Range("D26")
$r.NumberFormat = "@"
$r.Value = '2.799'
$r.Style = "Normal"
$ws.Range("E26").Value = '  +5.81%  '

$r = $ws.Range("D27")
$r.NumberFormat = "@"
$r.Value = '20.43'
$r.Style = "Normal"
$ws.Range("E27").Value = '  +5.42%  '

$r = $ws.Range("D28")
$r.NumberFormat = "@"
$r.Value = '149.64'
$r.Style = "Normal"
$ws.Range("E28").Value = '  -1.98%  '

$r = $ws.Range("D29")
$r.NumberFormat = "@"
$r.Value = '132.80'
$r.Style = "Normal"
$ws.Range("E29").Value = '  +4.20%  '

$ws.Range("B30").Value = 'WrappedliquidstakedEther2.0'
$ws.Range("C30").Value = 'https://coinranking.com/coin/CiixT63n3+wrappedliquidstakedether20-wsteth'
$ws.Range("D30").Value = '1.898.47'
$ws.Range("E30").Value = '  +1.65%  '

$ws.Range("B31").Value = 'ImmutableX'
$ws.Range("C31").Value = 'https://coinranking.com/coin/Z96jIvLU7+immutablex-imx'
$r = $ws.Range("D31")
$r.NumberFormat = "@"
$r.Value = '1.252'
$r.Style = "Normal"
$ws.Range("E31").Value = '  +29.87%  '

$r = $ws.Range("D32")
$r.NumberFormat = "@"
$r.Value = '6.826'
$r.Style = "Normal"
$ws.Range("E32").Value = '  +8.85%  '

$r = $ws.Range("D33")
$r.NumberFormat = "@"
$r.Value = '4.221'
$r.Style = "Normal"
$ws.Range("E33").Value = '  +3.88%  '

$r = $ws.Range("D34")
$r.NumberFormat = "@"
$r.Value = '13.92'
$r.Style = "Normal"
$ws.Range("E34").Value = '  +13.46%  '

$r = $ws.Range("D35")
$r.NumberFormat = "@"
$r.Value = '0.08776'
$r.Style = "Normal"
$ws.Range("E35").Value = '  +4.27%  '

$ws.Range("E36").Value = '  +4.52%  '

$r = $ws.Range("D37")
$r.NumberFormat = "@"
$r.Value = '5.637'
$r.Style = "Normal"
$ws.Range("E37").Value = '  +6.28%  '

$r = $ws.Range("D38")
$r.NumberFormat = "@"
$r.Value = '0.06653'
$r.Style = "Normal"
$ws.Range("E38").Value = '  +3.64%  '

$r = $ws.Range("D39")
$r.NumberFormat = "@"
$r.Value = '9.172'
$r.Style = "Normal"
$ws.Range("E39").Value = '  +3.63%  '

$r = $ws.Range("D40")
$r.NumberFormat = "@"
$r.Value = '0.02413'
$r.Style = "Normal"
$ws.Range("E40").Value = '  +4.35%  '

$r = $ws.Range("D41")
$r.NumberFormat = "@"
$r.Value = '0.2241'
$r.Style = "Normal"
$ws.Range("E41").Value = '  +7.24%  '

$r = $ws.Range("D42")
$r.NumberFormat = "@"
$r.Value = '1.277'
$r.Style = "Normal"
$ws.Range("E42").Value = '  +2.54%  '

$r = $ws.Range("D43")
$r.NumberFormat = "@"
$r.Value = '0.6468'
$r.Style = "Normal"
$ws.Range("E43").Value = '  +5.70%  '

$r = $ws.Range("D44")
$r.NumberFormat = "@"
$r.Value = '0.9983'
$r.Style = "Normal"
$ws.Range("E44").Value = '  -0.25%  '

$r = $ws.Range("D45")
$r.NumberFormat = "@"
$r.Value = '13.94'
$r.Style = "Normal"
$ws.Range("E45").Value = '  +6.35%  '

$r = $ws.Range("D46")
$r.NumberFormat = "@"
$r.Value = '0.6162'
$r.Style = "Normal"
$ws.Range("E46").Value = '  +4.36%  '

$r = $ws.Range("D47")
$r.NumberFormat = "@"
$r.Value = '3.839'
$r.Style = "Normal"
$ws.Range("E47").Value = '  +1.96%  '

$r = $ws.Range("D48")
$r.NumberFormat = "@"
$r.Value = '2.128'
$r.Style = "Normal"
$ws.Range("E48").Value = '  +5.59%  '

$r = $ws.Range("D49")
$r.NumberFormat = "@"
$r.Value = '129.67'
$r.Style = "Normal"
$ws.Range("E49").Value = '  +2.19%  '

$r = $ws.Range("D50")
$r.NumberFormat = "@"
$r.Value = '0.07326'
$r.Style = "Normal"
$ws.Range("E50").Value = '  +2.19%  '

$r = $ws.Range("D51")
$r.NumberFormat = "@"
$r.Value = '80.18'
$r.Style = "Normal"
$ws.Range("E51").Value = '  +6.12%  '
